# Update cryptocurrency price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.957.14"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "3.831.94"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "702.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").Value = "3.827.45"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "4.479.39"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "3.875.00"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "71.027.53"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("E22").Value = "  -4.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.735"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.177"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "3.794.58"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("E45").Value = "  +3.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "431.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.42%  "
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("E51").Value = "  -2.07%  "
